$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: J9 - drop the "FV: 4:45 AM BALLARD" line ---
$ws.Range("J9").Value = "M: 4:15 AM MEET COLLEGE (SOUTHWEST)`nIL: 5:00 AM AT IL OFFICE "

# --- Row 51: F51 - add IL / FV meet lines ---
$ws.Range("F51").Value = "M: 4:15 AM MEET SOUTH RACINE COURT (HWY 43 & Y)`nIL: 5:00 MEET AT IL OFFICE`nFV: 12:00 MEET"

# --- Row 64 / 67: G - driver / van note ---
$ws.Range("G64").Value = "Driver, 1/2`nSilver Van"
$ws.Range("G67").Value = "Driver, 1/2`nSilver Van"

# --- Rows 70-75: fill in F column details ---
$ws.Range("F70").Value = "6:00 AM START"
$ws.Range("F71").Value = "DC5-ITEM LEVEL"
$ws.Range("F72").Value = "AURORA OUTPATIENT RX #1300 WAUTOMA"
$ws.Range("F73").Value = "900 E DIVISION ST"
$ws.Range("F74").Value = "https://goo.gl/maps/JhUR6bvog6YmBrm48"
$ws.Range("F75").Value = "NEED SCANNER HOODS FOR THIS RX"

# --- Row 77: staffing note ---
$ws.Range("E77").Value = "1)"
$ws.Range("F77").Value = "Sarah"
$ws.Range("G77").Value = "@ Store. Equip, `nhelp at Pig after"

# --- Row 78: staffing note ---
$ws.Range("E78").Value = "2)"
$ws.Range("F78").Value = "Lori"
$ws.Range("G78").Value = "@ Store. help at Pig after"

# --- Row 81 ---
$ws.Range("E81").Value = " "
$ws.Range("F81").Value = "7:00 AM START"

# --- Rows 82-85 ---
$ws.Range("F82").Value = "DC5-FINANCIAL"
$ws.Range("F83").Value = "PIGGLY WIGGLY #311, WAUPACA"
$ws.Range("F84").Value = "810 W FULTON"
$ws.Range("F85").Value = "https://goo.gl/maps/vay5nnowsSN2"

# --- Row 87 ---
$ws.Range("E87").Value = "1)"
$ws.Range("F87").Value = "Jerry S"
$ws.Range("G87").Value = "@ Store, Equip"

# --- Row 88 ---
$ws.Range("E88").Value = "2)"
$ws.Range("F88").Value = "Lori"
$ws.Range("G88").Value = "@ Store, help after Aurora"

# --- New rows 89-95: additional staff list continuation ---
$ws.Range("E89").Value = "3)"
$ws.Range("F89").Value = "Sarah"
$ws.Range("G89").Value = "@ Store, help after Aurora"

$ws.Range("E90").Value = "4)"
$ws.Range("F90").Value = "Heather"
$ws.Range("G90").Value = "@ Store"

$ws.Range("E91").Value = "5)"
$ws.Range("F91").Value = "Katie"
$ws.Range("G91").Value = "@ Store, work w/ Serena"

$ws.Range("E92").Value = "6)"
$ws.Range("F92").Value = "Kirsten"
$ws.Range("G92").Value = "@ Store"

$ws.Range("E93").Value = "7)"
$ws.Range("F93").Value = "Marcia"
$ws.Range("G93").Value = "@ Store"

$ws.Range("E94").Value = "8)"
$ws.Range("F94").Value = "Michelle"
$ws.Range("G94").Value = "@ Store"

$ws.Range("E95").Value = "9)"
$ws.Range("F95").Value = "Serena"
$ws.Range("G95").Value = "@ Store, 1st Day, work w/ Katie"

# --- Rows 96-97 are blank spacer rows (nothing to set) ---

# --- Rows 98-103: new job block ---
$ws.Range("F98").Value = "6:30 AM START"
$ws.Range("F99").Value = "DC5-ITEM LEVEL"
$ws.Range("F100").Value = "KELLEY #63, SUN PRAIRIE MOBIL"
$ws.Range("F101").Value = "1010 DAVISON DR"
$ws.Range("F102").Value = "https://goo.gl/maps/AcJV9qEV2xQ2"
$ws.Range("F103").Value = "*IL Meet is 5:15 am at IL Office"

# --- Row 104 is a blank spacer row ---

# --- Rows 105-109: staffing list ---
$ws.Range("E105").Value = "1)"
$ws.Range("F105").Value = "Qiana"
$ws.Range("G105").Value = "@ Store,`nCamry, Equip"

$ws.Range("E106").Value = "2)"
$ws.Range("F106").Value = "Eva"

$ws.Range("E107").Value = "3)"
$ws.Range("F107").Value = "Evelin"

$ws.Range("E108").Value = "4)"
$ws.Range("F108").Value = "Josie"
$ws.Range("G108").Value = "@ Store"

$ws.Range("E109").Value = "5)"
$ws.Range("F109").Value = "Nate"
$ws.Range("G109").Value = "Driver, Optima"

# --- Rows 110-111 are blank spacer rows ---

# --- Rows 112-117: new job block ---
$ws.Range("F112").Value = "6:30 AM START"
$ws.Range("F113").Value = "MODAS"
$ws.Range("F114").Value = "KELLEY #58, TYLER CREEK MOBIL, HAMPSHIRE"
$ws.Range("F115").Value = "15N341 RTE 47"
$ws.Range("F116").Value = "https://goo.gl/maps/SW6S8XqhdGcLz6W2A"
$ws.Range("F117").Value = "*IL Meet is 5:45 am at IL Office"

# --- Row 118 is a blank spacer row ---

# --- Rows 119-124: staffing list ---
$ws.Range("E119").Value = "1)"
$ws.Range("F119").Value = "Mike G"
$ws.Range("G119").Value = "@ Store"

$ws.Range("E120").Value = "2)"
$ws.Range("F120").Value = "Angela"

$ws.Range("E121").Value = "3)"
$ws.Range("F121").Value = "Emily L"

$ws.Range("E122").Value = "4)"
$ws.Range("F122").Value = "Justin"
$ws.Range("G122").Value = "Driver,`nGray Van"

$ws.Range("E123").Value = "5)"
$ws.Range("F123").Value = "Krystin"

$ws.Range("E124").Value = "6)"
$ws.Range("F124").Value = "Taylor"
